# Update Name of Algo
# Apply updated values produced by the RandomForest imputation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value  = -7.369700000000002
$ws.Range("C9").Value  = -10.3096
$ws.Range("D12").Value = -7.192500000000006
$ws.Range("E15").Value = 16.07790000000001
$ws.Range("C18").Value = -12.53849999999999
$ws.Range("C20").Value = -11.78520000000001
$ws.Range("D26").Value = -8.326500000000003
$ws.Range("C27").Value = -12.5802
$ws.Range("D27").Value = -8.697500000000002
$ws.Range("D29").Value = -7.275899999999998
$ws.Range("D37").Value = -7.862299999999997
$ws.Range("D38").Value = -7.533300000000001
$ws.Range("E38").Value = 16.52299999999999
$ws.Range("E44").Value = 16.3347
$ws.Range("D51").Value = -8.559400000000002
$ws.Range("E51").Value = 16.4164
$ws.Range("D55").Value = -8.234499999999995
$ws.Range("E57").Value = 16.6647
$ws.Range("E63").Value = 18.23230000000002
$ws.Range("C69").Value = -11.3501
$ws.Range("D69").Value = -7.315299999999994
$ws.Range("D70").Value = -7.801400000000005
$ws.Range("E70").Value = 17.01730000000001
$ws.Range("C76").Value = -12.43410000000001
$ws.Range("C82").Value = -11.93599999999999
$ws.Range("D83").Value = -8.7501
$ws.Range("E99").Value = 16.4501
$ws.Range("D102").Value = -7.880300000000002
